$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1: update criterion weights ---
$ws.Range("D1").Value = 2
$ws.Range("F1").Value = 1
$ws.Range("H1").Value = 2

# --- Column I bonus comments & highlight colors ---

# Row 5: "thuật convex hull lạ" - change highlight from green to yellow
$ws.Range("I5").Interior.Color = 65535

# Row 6: "Nhóm tích cực làm bài" -> add "+ code OOP đẹp", change highlight yellow -> orange,
# taller row, and lower total score 10 -> 9.5
$ws.Range("I6").Value = "Nhóm tích cực làm bài + code OOP đẹp"
$ws.Range("I6").Interior.Color = 49407
$ws.Rows.Item(6).RowHeight = 30
$ws.Range("J6").Value = 9.5

# Row 8: "Nhóm tích cực tìm hiểu" -> add ", nhóm nộp sớm", taller row (highlight color unchanged)
$ws.Range("I8").Value = "Nhóm tích cực tìm hiểu, nhóm nộp sớm"
$ws.Rows.Item(8).RowHeight = 30

# Row 11: "Nhóm code đẹp, trực quan, thuật toán hướng tâm khá lạ" - green to yellow highlight
$ws.Range("I11").Interior.Color = 65535

# Row 13: "Nhóm tích cực tìm hiểu, nhóm nộp đầu tiên" - green to yellow highlight
$ws.Range("I13").Interior.Color = 65535

# --- Selection moves to I17 ---
$ws.Activate()
$ws.Range("I17").Select()
